# Update csv: added more authors, books and quantity set to 20k
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Set every existing "in_storage" quantity (column D, rows 2-55)
#    to 20000.
# ------------------------------------------------------------------
$ws.Range("D2:D55").Value = 20000

# ------------------------------------------------------------------
# 2. Append four new book rows (56-59), mirroring the style already
#    used by the existing data rows:
#      - column A uses the existing "s=2" numeric style elsewhere,
#        but the new trailing rows get a distinct boxed style,
#      - column B / C are shared-string lookups (cover type / price),
#      - column D is the quantity (20000 like every other row now),
#      - column E is the ISBN-like running counter.
# ------------------------------------------------------------------
$newRows = @(
    @{ Row = 56; A = 8893517191; B = "hardcover"; C = "11.80"; D = 20000; E = 1234567891068 },
    @{ Row = 57; A = 8893517192; B = "paperback"; C = "20.37"; D = 20000; E = 1234567891069 },
    @{ Row = 58; A = 8893517192; B = "hardcover"; C = "31.15"; D = 20000; E = 1234567891070 },
    @{ Row = 59; A = 8893517193; B = "paperback"; C = "12.30"; D = 20000; E = 1234567891071 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E

    # Column C keeps the same right-aligned style already used for the
    # price column on every other row.
    $ws.Cells.Item($row, 3).Style = $ws.Cells.Item(53, 3).Style

    # Column E keeps the existing ISBN column style.
    $ws.Cells.Item($row, 5).Style = $ws.Cells.Item(53, 5).Style

    # Column A gets the new boxed / right-aligned / wrapped look.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Font.Name = "Arial"
    $cellA.Font.Size = 10
    $cellA.HorizontalAlignment = -4152
    $cellA.WrapText = $true
    $cellA.Borders.LineStyle = 1
    $cellA.Borders.Weight = -4138
    $cellA.Borders.Color = 13421772

    $ws.Rows.Item($row).RowHeight = 15
}

# ------------------------------------------------------------------
# 3. Refresh the window's selection/scroll position to match where the
#    author ended up after appending the new rows.
# ------------------------------------------------------------------
$ws.Range("G58").Select()
$excel.ActiveWindow.ScrollRow = 38
